$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1563988.2
$ws.Range("I132").Value = 1424.7451
$ws.Range("K132").Value = 4274.2353
$ws.Range("M132").Value = -1744.2353

$ws.Range("H135").Value = 1231.5385
$ws.Range("I135").Value = 1287.3182
$ws.Range("J135").Value = 924.75
$ws.Range("K135").Value = 11585.8638
$ws.Range("L135").Value = 8322.75
$ws.Range("M135").Value = -9050.863799999999
$ws.Range("N135").Value = -13392.75

$ws.Range("H137").Value = 1034.1957
$ws.Range("I137").Value = 897.025
$ws.Range("J137").Value = 1948.6666
$ws.Range("K137").Value = 2691.075
$ws.Range("L137").Value = 5845.9998
$ws.Range("M137").Value = -141.0749999999998
$ws.Range("N137").Value = -10945.9998

$ws.Range("H138").Value = 3760.307
$ws.Range("I138").Value = 1772.0625
$ws.Range("J138").Value = 4896.4463
$ws.Range("K138").Value = 5316.1875
$ws.Range("L138").Value = 14689.3389
$ws.Range("M138").Value = -176.1875
$ws.Range("N138").Value = -24969.3389

$ws.Range("H141").Value = 578.75
$ws.Range("I141").Value = 578.75
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 1736.25
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3443.75
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13137.37
$ws.Range("I32").Value = 11422.67
$ws.Range("J32").Value = 40001
$ws.Range("K32").Value = 11422.67
$ws.Range("L32").Value = 40001
$ws.Range("M32").Value = -11135.67
$ws.Range("N32").Value = -40575

$ws.Range("H102").Value = 7420
$ws.Range("I102").Value = 2904
$ws.Range("J102").Value = 30000
$ws.Range("K102").Value = 2904
$ws.Range("L102").Value = 30000
$ws.Range("M102").Value = -1282
$ws.Range("N102").Value = -33244

$ws.Range("H132").Value = 1566.9803
$ws.Range("I132").Value = 1135.5
$ws.Range("J132").Value = 2602.5334
$ws.Range("K132").Value = 3406.5
$ws.Range("L132").Value = 7807.600199999999
$ws.Range("M132").Value = -876.5
$ws.Range("N132").Value = -12867.6002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 50003190
$ws.Range("I99").Value = 58826520
$ws.Range("J99").Value = 4327
$ws.Range("K99").Value = 58826520
$ws.Range("L99").Value = 4327
$ws.Range("M99").Value = -58825022
$ws.Range("N99").Value = -7323

$ws.Range("H107").Value = 4917.5483
$ws.Range("I107").Value = 628.61536
$ws.Range("J107").Value = 27220
$ws.Range("K107").Value = 628.61536
$ws.Range("L107").Value = 27220
$ws.Range("M107").Value = 1291.38464
$ws.Range("N107").Value = -31060

$ws.Range("H134").Value = 18095.328
$ws.Range("I134").Value = 1569.4108
$ws.Range("J134").Value = 203185.6
$ws.Range("K134").Value = 4708.232400000001
$ws.Range("L134").Value = 609556.8
$ws.Range("M134").Value = -2173.232400000001
$ws.Range("N134").Value = -614626.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2498.2078
$ws.Range("I31").Value = 1878.8545
$ws.Range("J31").Value = 4046.5908
$ws.Range("K31").Value = 1878.8545
$ws.Range("L31").Value = 4046.5908
$ws.Range("M31").Value = -1583.8545
$ws.Range("N31").Value = -4636.5908

$ws.Range("H34").Value = 2498.2078
$ws.Range("I34").Value = 1878.8545
$ws.Range("J34").Value = 4046.5908
$ws.Range("K34").Value = 1878.8545
$ws.Range("L34").Value = 4046.5908
$ws.Range("M34").Value = -1676.8545
$ws.Range("N34").Value = -4450.5908

$ws.Range("H99").Value = 2178.1936
$ws.Range("I99").Value = 2150.6
$ws.Range("J99").Value = 2228.3635
$ws.Range("K99").Value = 2150.6
$ws.Range("L99").Value = 2228.3635
$ws.Range("M99").Value = -652.5999999999999
$ws.Range("N99").Value = -5224.363499999999

$ws.Range("H126").Value = 2178.1936
$ws.Range("I126").Value = 2150.6
$ws.Range("J126").Value = 2228.3635
$ws.Range("K126").Value = 6451.799999999999
$ws.Range("L126").Value = 6685.0905
$ws.Range("M126").Value = -3981.799999999999
$ws.Range("N126").Value = -11625.0905

$ws.Range("H138").Value = 77705
$ws.Range("J138").Value = 77705
$ws.Range("L138").Value = 77705
$ws.Range("N138").Value = -87985

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 466.66666
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 1200
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = 48
$ws.Range("N92").Value = -3996

$ws.Range("H107").Value = 289056.34
$ws.Range("J107").Value = 432868.84
$ws.Range("L107").Value = 1298606.52
$ws.Range("N107").Value = -1302446.52

$ws.Range("H129").Value = 43379.4
$ws.Range("I129").Value = 1391.1111
$ws.Range("J129").Value = 66997.81
$ws.Range("K129").Value = 4173.3333
$ws.Range("L129").Value = 200993.43
$ws.Range("M129").Value = 826.6666999999998
$ws.Range("N129").Value = -210993.43

$ws.Range("H131").Value = 41374.7
$ws.Range("I131").Value = 202406
$ws.Range("J131").Value = 23482.334
$ws.Range("K131").Value = 607218
$ws.Range("L131").Value = 70447.00199999999
$ws.Range("M131").Value = -602178
$ws.Range("N131").Value = -80527.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11972.94
$ws.Range("I70").Value = 16371.5
$ws.Range("J70").Value = 4153.278
$ws.Range("K70").Value = 16371.5
$ws.Range("L70").Value = 4153.278
$ws.Range("M70").Value = -16101.5
$ws.Range("N70").Value = -4693.278

$ws.Range("H73").Value = 11972.94
$ws.Range("I73").Value = 16371.5
$ws.Range("J73").Value = 4153.278
$ws.Range("K73").Value = 16371.5
$ws.Range("L73").Value = 4153.278
$ws.Range("M73").Value = -15435.5
$ws.Range("N73").Value = -6025.278

$ws.Range("H126").Value = 2527.577
$ws.Range("I126").Value = 1905.3182
$ws.Range("J126").Value = 5950
$ws.Range("K126").Value = 5715.9546
$ws.Range("L126").Value = 17850
$ws.Range("M126").Value = -3245.9546
$ws.Range("N126").Value = -22790

$ws.Range("H132").Value = 2312.3408
$ws.Range("I132").Value = 1916.7778
$ws.Range("J132").Value = 4092.375
$ws.Range("K132").Value = 5750.3334
$ws.Range("L132").Value = 12277.125
$ws.Range("M132").Value = -3220.3334
$ws.Range("N132").Value = -17337.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3205.875
$ws.Range("I136").Value = 654.1852
$ws.Range("J136").Value = 8505.538
$ws.Range("K136").Value = 1962.5556
$ws.Range("L136").Value = 25516.614
$ws.Range("M136").Value = 587.4443999999999
$ws.Range("N136").Value = -30616.614
